# The edit reorders the 21 data rows (rows 2-15 and 17-23; the header row 1
# and row 16 are untouched) onto a different set of dates/values — i.e. the
# content of each data row is replaced by the content that used to live in
# a different row. Snapshot every source row's full A:T values first (since
# the reassignment contains cycles), then write them back in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# newRow -> oldRow : row "newRow" ends up holding what used to be in "oldRow"
$mapping = @{
    2  = 21
    3  = 22
    4  = 7
    5  = 2
    6  = 3
    7  = 4
    8  = 19
    9  = 20
    10 = 17
    11 = 18
    12 = 9
    13 = 10
    14 = 6
    15 = 23
    17 = 11
    18 = 12
    19 = 13
    20 = 14
    21 = 15
    22 = 5
    23 = 8
}

# Snapshot all source rows (A:T) before writing anything, since the mapping
# contains cycles and naive in-place writes would clobber data we still need.
$snapshot = @{}
foreach ($oldRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($oldRow)) {
        $snapshot[$oldRow] = $ws.Range("A$oldRow`:T$oldRow").Value2
    }
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Range("A$newRow`:T$newRow").Value2 = $snapshot[$oldRow]
}
